$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D = Price column, E = Volume(1h) column.
# Values that look like plain numbers get a leading apostrophe so Excel
# stores them as text (preserving exact formatting, e.g. trailing zeros),
# matching the inline-string cells in the source data.

$ws.Cells.Item(2, 4).Value = "26.167.00"
$ws.Cells.Item(2, 5).Value = "  -4.18%  "
$ws.Cells.Item(3, 4).Value = "1.659.69"
$ws.Cells.Item(4, 4).Value = "'1.005"
$ws.Cells.Item(4, 5).Value = "  +0.08%  "
$ws.Cells.Item(5, 4).Value = "'218.42"
$ws.Cells.Item(5, 5).Value = "  -2.40%  "
$ws.Cells.Item(6, 4).Value = "'0.5164"
$ws.Cells.Item(6, 5).Value = "  -2.65%  "
$ws.Cells.Item(7, 4).Value = "'1.005"
$ws.Cells.Item(7, 5).Value = "  +0.18%  "
$ws.Cells.Item(8, 4).Value = "'0.2591"
$ws.Cells.Item(8, 5).Value = "  -2.38%  "
$ws.Cells.Item(9, 5).Value = "  -1.79%  "
$ws.Cells.Item(10, 4).Value = "'19.95"
$ws.Cells.Item(10, 5).Value = "  -3.75%  "
$ws.Cells.Item(11, 4).Value = "'0.07792"
$ws.Cells.Item(11, 5).Value = "  +2.11%  "
$ws.Cells.Item(12, 4).Value = "1.660.18"
$ws.Cells.Item(12, 5).Value = "  -2.57%  "
$ws.Cells.Item(13, 4).Value = "'4.298"
$ws.Cells.Item(13, 5).Value = "  -4.75%  "
$ws.Cells.Item(14, 4).Value = "1.887.37"
$ws.Cells.Item(14, 5).Value = "  -2.69%  "
$ws.Cells.Item(15, 4).Value = "'0.5560"
$ws.Cells.Item(15, 5).Value = "  -3.62%  "
$ws.Cells.Item(16, 4).Value = "0.0₅8076"
$ws.Cells.Item(16, 5).Value = "  -0.62%  "
$ws.Cells.Item(17, 4).Value = "'64.28"
$ws.Cells.Item(17, 5).Value = "  -4.82%  "
$ws.Cells.Item(18, 4).Value = "26.205.26"
$ws.Cells.Item(18, 5).Value = "  -4.03%  "
$ws.Cells.Item(19, 4).Value = "'212.13"
$ws.Cells.Item(19, 5).Value = "  -1.36%  "
$ws.Cells.Item(20, 5).Value = "  +0.23%  "
$ws.Cells.Item(21, 4).Value = "'4.421"
$ws.Cells.Item(21, 5).Value = "  -4.16%  "
$ws.Cells.Item(22, 4).Value = "'10.05"
$ws.Cells.Item(22, 5).Value = "  -2.93%  "
$ws.Cells.Item(23, 4).Value = "'5.967"
$ws.Cells.Item(23, 5).Value = "  +0.23%  "
$ws.Cells.Item(24, 5).Value = "  +0.19%  "
$ws.Cells.Item(25, 4).Value = "'144.36"
$ws.Cells.Item(25, 5).Value = "  +0.52%  "
$ws.Cells.Item(26, 4).Value = "'1.756"
$ws.Cells.Item(26, 5).Value = "  +2.38%  "
$ws.Cells.Item(27, 4).Value = "'0.1165"
$ws.Cells.Item(27, 5).Value = "  -3.03%  "
$ws.Cells.Item(28, 4).Value = "'6.980"
$ws.Cells.Item(28, 5).Value = "  -3.15%  "
$ws.Cells.Item(29, 4).Value = "'15.84"
$ws.Cells.Item(29, 5).Value = "  -1.45%  "
$ws.Cells.Item(30, 4).Value = "'0.05275"
$ws.Cells.Item(30, 5).Value = "  -1.74%  "
$ws.Cells.Item(31, 4).Value = "'1.256"
$ws.Cells.Item(31, 5).Value = "  -2.45%  "
$ws.Cells.Item(32, 4).Value = "'3.366"
$ws.Cells.Item(32, 5).Value = "  -2.82%  "
$ws.Cells.Item(33, 4).Value = "'3.221"
$ws.Cells.Item(33, 5).Value = "  -5.25%  "
$ws.Cells.Item(34, 4).Value = "'1.571"
$ws.Cells.Item(34, 5).Value = "  -4.00%  "
$ws.Cells.Item(35, 5).Value = "  -3.74%  "
$ws.Cells.Item(36, 5).Value = "  -1.91%  "
$ws.Cells.Item(37, 4).Value = "'0.9279"
$ws.Cells.Item(37, 5).Value = "  -1.78%  "
$ws.Cells.Item(38, 4).Value = "1.169.66"
$ws.Cells.Item(38, 5).Value = "  +12.61%  "
$ws.Cells.Item(39, 4).Value = "'0.5666"
$ws.Cells.Item(39, 5).Value = "  -2.06%  "
$ws.Cells.Item(40, 5).Value = "  -2.02%  "
$ws.Cells.Item(41, 4).Value = "'1.005"
$ws.Cells.Item(41, 5).Value = "  +0.17%  "
$ws.Cells.Item(42, 4).Value = "'0.8443"
$ws.Cells.Item(43, 4).Value = "'5.694"
$ws.Cells.Item(43, 5).Value = "  -1.06%  "
$ws.Cells.Item(44, 4).Value = "'100.53"
$ws.Cells.Item(44, 5).Value = "  -0.37%  "
$ws.Cells.Item(45, 4).Value = "1.797.34"
$ws.Cells.Item(45, 5).Value = "  -2.75%  "
$ws.Cells.Item(46, 4).Value = "0.0₈113"
$ws.Cells.Item(46, 5).Value = "  -3.13%  "
$ws.Cells.Item(47, 5).Value = "  +0.41%  "
$ws.Cells.Item(48, 4).Value = "'55.89"
$ws.Cells.Item(48, 5).Value = "  -3.03%  "
$ws.Cells.Item(49, 5).Value = "  +0.28%  "
$ws.Cells.Item(50, 4).Value = "'7.889"
$ws.Cells.Item(50, 5).Value = "  -2.22%  "
$ws.Cells.Item(51, 4).Value = "'0.05055"
$ws.Cells.Item(51, 5).Value = "  -3.36%  "
